# Update the weekly data rows: columns D (Fecha), J (Volumen), K (Precio minimo),
# L (Precio maximo), M (Precio promedio ponderado) and P (Precio $/Kg) are
# rotated across the data rows (2-9). Other columns remain unchanged.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @{ D = 44229; J = 120; K = 44000; L = 45000; M = 44500; P = 3423 }
    3 = @{ D = 44159; J = 100; K = 23000; L = 24000; M = 23500; P = 1808 }
    4 = @{ D = 44397; J = 140; K = 12500; L = 13000; M = 12750; P = 981 }
    5 = @{ D = 44406; J = 160; K = 17000; L = 18000; M = 17500; P = 1346 }
    6 = @{ D = 44320; J = 160; K = 19000; L = 20000; M = 19500; P = 1500 }
    7 = @{ D = 44469; J = 140; K = 13000; L = 14000; M = 13500; P = 1038 }
    8 = @{ D = 44379; J = 120; K = 12000; L = 13000; M = 12667; P = 974 }
    9 = @{ D = 44389; J = 120; K = 12000; L = 13000; M = 12500; P = 962 }
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Range("D$row").Value = $vals.D
    $ws.Range("J$row").Value = $vals.J
    $ws.Range("K$row").Value = $vals.K
    $ws.Range("L$row").Value = $vals.L
    $ws.Range("M$row").Value = $vals.M
    $ws.Range("P$row").Value = $vals.P
}
